$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, pushing the current rows 7-8 down to 8-9.
$ws.Rows.Item(7).Insert()

# Copy the date-format style (used by column D throughout) onto the new D7 cell.
$ws.Range("D8").Copy()
$ws.Range("D7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Populate the new row 7 with the new record's values.
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Macroferia Regional de Talca"
$ws.Range("C7").Value = "Maule"
$ws.Range("D7").Value = 45264
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100104
$ws.Range("H7").Value = "Frutos de pepita"
$ws.Range("I7").Value = 100104004
$ws.Range("J7").Value = "Níspero"
$ws.Range("K7").Value = "Golden Nugget"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 150
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("Q7").Value = "`$/caja 10 kilos"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1500
$ws.Range("T7").Value = 10

# Row 8 now holds what used to be row 7's data, except the variety/volume swap
# to Californiana(o)/100 (the diff shows row 8 ending up with the old row 7's
# other fields but the variety+volume that used to belong to row 8).
$ws.Range("K8").Value = "Californiana(o)"
$ws.Range("M8").Value = 100

# Row 9 (shifted from the old row 8) keeps its original values already, as
# carried down by the row insert/shift - nothing further required there.
